# Daily 100 Error Counts - append newly reported days (33-36) and
# move the sheet's viewport/selection down to the new bottom row,
# matching the "Add files via upload" refresh of the dashboard data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily rows: Date, Total Count, Session Timeout Errors, Errors Requiring Analysis
$newRows = @(
    @{ Row = 33; Date = 45982; Total = 539; Timeout = 19; Analysis = 520 },
    @{ Row = 34; Date = 45985; Total = 653; Timeout = 34; Analysis = 619 },
    @{ Row = 35; Date = 45986; Total = 672; Timeout = 17; Analysis = 655 },
    @{ Row = 36; Date = 45987; Total = 581; Timeout = 17; Analysis = 564 }
)

foreach ($r in $newRows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.Date
    $ws.Cells.Item($r.Row, 2).Value = $r.Total
    $ws.Cells.Item($r.Row, 3).Value = $r.Timeout
    $ws.Cells.Item($r.Row, 4).Value = $r.Analysis
}

# Scroll / reselect to reflect the new bottom-of-data row, as Excel does
# when the user scrolls down after adding rows and saves.
$excel.ActiveWindow.ScrollRow = 29
[void]$ws.Range("A36:D36").Select()
